$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1676.0769
$ws.Range("I19").Value = 1396.3334
$ws.Range("J19").Value = 1760
$ws.Range("K19").Value = 1396.3334
$ws.Range("L19").Value = 1760
$ws.Range("M19").Value = -1221.3334
$ws.Range("N19").Value = -2110
$ws.Range("H32").Value = 800
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H64").Value = 145742.86
$ws.Range("I64").Value = 252250
$ws.Range("J64").Value = 3733.3333
$ws.Range("K64").Value = 252250
$ws.Range("L64").Value = 3733.3333
$ws.Range("M64").Value = -252002
$ws.Range("N64").Value = -4229.3333
$ws.Range("H67").Value = 145742.86
$ws.Range("I67").Value = 252250
$ws.Range("J67").Value = 3733.3333
$ws.Range("K67").Value = 252250
$ws.Range("L67").Value = 3733.3333
$ws.Range("M67").Value = -251392
$ws.Range("N67").Value = -5449.3333
$ws.Range("H70").Value = 1045.091
$ws.Range("I70").Value = 1074.625
$ws.Range("J70").Value = 1017.2941
$ws.Range("K70").Value = 3223.875
$ws.Range("L70").Value = 3051.8823
$ws.Range("M70").Value = -2953.875
$ws.Range("N70").Value = -3591.8823
$ws.Range("H73").Value = 1045.091
$ws.Range("I73").Value = 1074.625
$ws.Range("J73").Value = 1017.2941
$ws.Range("K73").Value = 3223.875
$ws.Range("L73").Value = 3051.8823
$ws.Range("M73").Value = -2287.875
$ws.Range("N73").Value = -4923.882299999999
$ws.Range("H74").Value = 3284.7144
$ws.Range("I74").Value = 3150.75
$ws.Range("J74").Value = 3463.3333
$ws.Range("K74").Value = 3150.75
$ws.Range("L74").Value = 3463.3333
$ws.Range("M74").Value = -2214.75
$ws.Range("N74").Value = -5335.3333
$ws.Range("H77").Value = 3284.7144
$ws.Range("I77").Value = 3150.75
$ws.Range("J77").Value = 3463.3333
$ws.Range("K77").Value = 15753.75
$ws.Range("L77").Value = 17316.6665
$ws.Range("M77").Value = -11073.75
$ws.Range("N77").Value = -26676.6665
$ws.Range("H88").Value = 6287.5
$ws.Range("I88").Value = 3750
$ws.Range("J88").Value = 7133.3335
$ws.Range("K88").Value = 3750
$ws.Range("L88").Value = 7133.3335
$ws.Range("M88").Value = -3344
$ws.Range("N88").Value = -7945.3335
$ws.Range("H91").Value = 6287.5
$ws.Range("I91").Value = 3750
$ws.Range("J91").Value = 7133.3335
$ws.Range("K91").Value = 3750
$ws.Range("L91").Value = 7133.3335
$ws.Range("M91").Value = -2346
$ws.Range("N91").Value = -9941.333500000001
$ws.Range("H107").Value = 616.4545000000001
$ws.Range("I107").Value = 616.4545000000001
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 616.4545000000001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1303.5455
$ws.Range("N107").ClearContents()
$ws.Range("H118").Value = 9160.833000000001
$ws.Range("I118").Value = 11714.444
$ws.Range("K118").Value = 35143.33199999999
$ws.Range("M118").Value = -33486.33199999999
$ws.Range("H132").Value = 5560988.5
$ws.Range("I132").Value = 6255484.5
$ws.Range("K132").Value = 18766453.5
$ws.Range("M132").Value = -18763923.5
$ws.Range("H138").Value = 2948.2195
$ws.Range("I138").Value = 2860.4443
$ws.Range("J138").Value = 2972.9062
$ws.Range("K138").Value = 8581.332900000001
$ws.Range("L138").Value = 8918.7186
$ws.Range("M138").Value = -3441.332900000001
$ws.Range("N138").Value = -19198.7186

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 47958.418
$ws.Range("I32").Value = 11896.909
$ws.Range("J32").Value = 166961.4
$ws.Range("K32").Value = 11896.909
$ws.Range("L32").Value = 166961.4
$ws.Range("M32").Value = -11609.909
$ws.Range("N32").Value = -167535.4
$ws.Range("H74").Value = 1106.0883
$ws.Range("I74").Value = 1097.75
$ws.Range("J74").Value = 1145
$ws.Range("K74").Value = 1097.75
$ws.Range("L74").Value = 1145
$ws.Range("M74").Value = -223.75
$ws.Range("N74").Value = -2893
$ws.Range("H77").Value = 1106.0883
$ws.Range("I77").Value = 1097.75
$ws.Range("J77").Value = 1145
$ws.Range("K77").Value = 5488.75
$ws.Range("L77").Value = 5725
$ws.Range("M77").Value = -1120.75
$ws.Range("N77").Value = -14461
$ws.Range("H88").Value = 1483.8334
$ws.Range("I88").Value = 1466.3334
$ws.Range("J88").Value = 1489.6666
$ws.Range("K88").Value = 1466.3334
$ws.Range("L88").Value = 1489.6666
$ws.Range("M88").Value = -1060.3334
$ws.Range("N88").Value = -2301.6666
$ws.Range("H91").Value = 1483.8334
$ws.Range("I91").Value = 1466.3334
$ws.Range("J91").Value = 1489.6666
$ws.Range("K91").Value = 1466.3334
$ws.Range("L91").Value = 1489.6666
$ws.Range("M91").Value = -62.33339999999998
$ws.Range("N91").Value = -4297.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1775.25
$ws.Range("I134").Value = 1438.5454
$ws.Range("J134").Value = 2516
$ws.Range("K134").Value = 4315.6362
$ws.Range("L134").Value = 7548
$ws.Range("M134").Value = -1780.6362
$ws.Range("N134").Value = -12618

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 512.75
$ws.Range("I25").Value = 125.5
$ws.Range("J25").Value = 900
$ws.Range("K25").Value = 376.5
$ws.Range("L25").Value = 2700
$ws.Range("M25").Value = -207.5
$ws.Range("N25").Value = -3038
$ws.Range("H30").Value = 512.75
$ws.Range("I30").Value = 125.5
$ws.Range("J30").Value = 900
$ws.Range("K30").Value = 376.5
$ws.Range("L30").Value = 2700
$ws.Range("M30").Value = -274.5
$ws.Range("N30").Value = -2904
$ws.Range("H113").Value = 800.0833
$ws.Range("I113").Value = 1165
$ws.Range("J113").Value = 539.4286
$ws.Range("K113").Value = 3495
$ws.Range("L113").Value = 1618.2858
$ws.Range("M113").Value = -1325
$ws.Range("N113").Value = -5958.2858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 72324.64
$ws.Range("I107").Value = 693.55554
$ws.Range("J107").Value = 201260.6
$ws.Range("K107").Value = 2080.66662
$ws.Range("L107").Value = 603781.8
$ws.Range("M107").Value = -160.66662
$ws.Range("N107").Value = -607621.8
$ws.Range("H113").Value = 354.11765
$ws.Range("I113").Value = 322.54544
$ws.Range("J113").Value = 412
$ws.Range("K113").Value = 967.63632
$ws.Range("L113").Value = 1236
$ws.Range("M113").Value = 1202.36368
$ws.Range("N113").Value = -5576
$ws.Range("H122").Value = 2150
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2150
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6450
$ws.Range("M122").ClearContents()
